$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("556-muzzles")

# Update R and U columns for rows 23-29
$ws.Range("R23").Value = -14
$ws.Range("U23").Value = -0.2

$ws.Range("R24").Value = -16
$ws.Range("U24").Value = -0.15

$ws.Range("R25").Value = -17
$ws.Range("U25").Value = -0.1

$ws.Range("R26").Value = -17
$ws.Range("U26").Value = -0.1

$ws.Range("R27").Value = -16
$ws.Range("U27").Value = -0.05

$ws.Range("R28").Value = -16
$ws.Range("U28").Value = -0.05

$ws.Range("R29").Value = -15
$ws.Range("U29").Value = -0.15

# Update the active selection to match the saved view state
$ws.Range("P27").Select()
